$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header text: MODEL_CONDITION -> MODELCONDITION (currently in E1)
$ws.Range("E1").Value = "MODELCONDITION"

# Delete entire column A (values 3,8,14,16,18 with border style), shifting remaining columns left
$ws.Columns("A").Delete()
